# Update F-column ("想去人数") values across the sheets that report
# event attendance counts: 展览 (sheet1), 演出 (sheet2) and the merged
# 全部类型 (sheet4) view. 本地生活 (sheet3) has no data rows so it is
# left untouched.

$wb = $excel.ActiveWorkbook

$sheet1Values = @{
    2 = 621
    3 = 206
    4 = 592
    5 = 538
    6 = 300
    7 = 2704
    8 = 461
    9 = 7521
    10 = 197
    11 = 463
    12 = 29
    13 = 256
    14 = 42
}

$sheet2Values = @{
    2 = 13
    3 = 18
    4 = 2
    5 = 1
}

$sheet4Values = @{
    2 = 621
    3 = 206
    4 = 592
    5 = 538
    6 = 300
    7 = 13
    8 = 18
    9 = 2704
    10 = 461
    11 = 7521
    12 = 197
    13 = 463
    14 = 29
    15 = 2
    16 = 1
    17 = 256
    18 = 42
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Values.Keys) {
    $ws1.Range("F$row").Value = $sheet1Values[$row]
}

$ws2 = $wb.Worksheets.Item("演出")
foreach ($row in $sheet2Values.Keys) {
    $ws2.Range("F$row").Value = $sheet2Values[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Values.Keys) {
    $ws4.Range("F$row").Value = $sheet4Values[$row]
}
